$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the five new rows needed to grow the test-case table from 20 to
#    25 data rows. Insert from the bottom up so earlier row numbers (used as
#    insertion anchors) stay valid.
# ---------------------------------------------------------------------------
$ws.Rows("21:21").Insert()   # becomes the split of old row 20 (navigate <>)
$ws.Rows("17:18").Insert()   # two brand-new "frame size" rows
$ws.Rows("16:16").Insert()   # split of the old "frames...size" row
$ws.Rows("7:7").Insert()     # split of the old "commas...price" row

# ---------------------------------------------------------------------------
# 2. Clone per-column formatting onto the freshly inserted (blank) rows from
#    a neighbouring row that already carries the right style.
# ---------------------------------------------------------------------------
$ws.Range("B6:J6").Copy()
$ws.Range("B7:J7").PasteSpecial(-4122)

$ws.Range("B15:J15").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)

$ws.Range("B16:J16").Copy()
$ws.Range("B17:J20").PasteSpecial(-4122)

$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Row heights — rows 6 & 7 (the two "commas" rows) grew to 40pt.
# ---------------------------------------------------------------------------
$ws.Rows("6:7").RowHeight = 40

# ---------------------------------------------------------------------------
# 4. Re-sequence the ID column (B6:B25) with the same "+1" formula pattern
#    the sheet already used further down.
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 1
$ws.Range("B7").Formula = "=B6+1"
$ws.Range("B8").Formula = "=B7+1"
$ws.Range("B9").Formula = "=B8+1"
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 6
$ws.Range("B12").Value = 7
$ws.Range("B13").Value = 8
$ws.Range("B15").Formula = "=B13+1"
$ws.Range("B16").Formula = "=B15+1"
$ws.Range("B17").Formula = "=B16+1"
$ws.Range("B18").Formula = "=B17+1"
$ws.Range("B19").Value = 13
$ws.Range("B20").Formula = "=B19+1"
$ws.Range("B21").Value = 15
$ws.Range("B22").Value = 16
$ws.Range("B23").Value = 17
$ws.Range("B24").Value = 18
$ws.Range("B25").Value = 19

# ---------------------------------------------------------------------------
# 5. Test-case description text (column C) for every data row.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = "Verify that commas are used as decimal separators to separate groups of thousands, millions, `nbillions, etc. for original price"
$ws.Range("C7").Value = "Verify that commas are used as decimal separators to separate groups of thousands, millions, `nbillions, etc. for discounted price"
$ws.Range("C8").Value = "Verify that discounted price is rounded to the nearest integer"
$ws.Range("C9").Value = "Verify that currency displayed is correct"
$ws.Range("C10").Value = "Verify that the text is in the correct font"
$ws.Range("C11").Value = "Verify that the text is in the correct size"
$ws.Range("C12").Value = "Verify that the text is in the correct colour"
$ws.Range("C13").Value = "Verify that the spelling is correct"
$ws.Range("C15").Value = "Verify that the photo list cannot display less than 1 photo"
$ws.Range("C16").Value = "Verify that the photo list can display from 1 to 5 photos"
$ws.Range("C17").Value = "Verify that the photo list cannot display more than 5 photos"
$ws.Range("C18").Value = "Verify that the first photo is displayed on the big frame"
$ws.Range("C19").Value = "Verify that the frame for small photos are in the correct size"
$ws.Range("C20").Value = "Verify that the frame for big photos are in the correct size"
$ws.Range("C21").Value = "Verify that the <> button is in the correct size"
$ws.Range("C22").Value = "Verify that the <> button is in the correct colour"
$ws.Range("C23").Value = "Verify that the <> button is properly aligned"
$ws.Range("C24").Value = "Verify that user can click on the < button to navigate to the previous photos"
$ws.Range("C25").Value = "Verify that user can click on the > button to navigate to the next photos"

# ---------------------------------------------------------------------------
# 6. Column / view cosmetics captured by the diff.
# ---------------------------------------------------------------------------
$ws.Range("C1").ColumnWidth = 49.81640625
$ws.Application.ActiveWindow.Zoom = 115
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C17").Select()
